$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.94498
$ws.Range("H2").Value = 17.83494
$ws.Range("I2").Value = 0.4679240463447598
$ws.Range("J2").Value = 0.4679240463447597
$ws.Range("Q2").Value = 0.03775260466
$ws.Range("R2").Value = 0.33977344194
$ws.Range("S2").Value = 0.4679240463447598
$ws.Range("T2").Value = 0.4679240463447597

# Row 3
$ws.Range("H3").Value = 8.352077
$ws.Range("I3").Value = 0.219128164447035
$ws.Range("J3").Value = 0.219128164447035
$ws.Range("S3").Value = 0.219128164447035
$ws.Range("T3").Value = 0.219128164447035

# Row 4
$ws.Range("G4").Value = 3.976005
$ws.Range("H4").Value = 11.928015
$ws.Range("I4").Value = 0.3129477892082053
$ws.Range("J4").Value = 0.3129477892082053
$ws.Range("Q4").Value = 0.025248957085
$ws.Range("R4").Value = 0.227240613765
$ws.Range("S4").Value = 0.3129477892082053
$ws.Range("T4").Value = 0.3129477892082053
